$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37:H37").Copy()
$ws.Range("A38:H38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A38").Value = "2025-08-21 03:50:50 UTC"
$ws.Range("B38").Value = "2025-08-21 09:20:50 IST"
$ws.Range("C38").Value = "UPDATED"
$ws.Range("D38").Value = "New circular processed."
$ws.Range("E38").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Range("F38").Value = "INGOT-21-08-2025.pdf"
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 3

Write-Host "Done"
